$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.149.43'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.90%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.895.46'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.81%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.21%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '485.14'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.34%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.98'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.73%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.622'
$ws.Range("D7").Style = "Normal"

# Row 8
$ws.Range("E8").Value = '  +0.03%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.736'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.27%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.175'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.42%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000354'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.52%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.99'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.62%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.64'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.20%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.525.04'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.07%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.899.23'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.75%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.33'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.83%  '

# Row 17
$ws.Range("E17").Value = '  -0.23%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.28'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.77%  '

# Row 19
$ws.Range("E19").Value = '  +2.11%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.177.35'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.82%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '429.55'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.86%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.56'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +8.27%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.90'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.14%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.98'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.55%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.43'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +14.36%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.71'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.98%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.03'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +9.54%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.50'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.19%  '

# Row 29
$ws.Range("E29").Value = '  -2.06%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '719.29'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.48%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.69'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.63%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.130'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.79%  '

# Row 33
$ws.Range("E33").Value = '  +4.64%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0900'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.51%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '41.46'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.21%  '

# Row 36
$ws.Range("E36").Value = '  +14.75%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '61.36'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.39%  '

# Row 38
$ws.Range("B38").Value = 'TheGraph'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.397'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +18.05%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.145'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.81%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.998'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.01%  '

# Row 41
$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.01'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +10.38%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0496'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +7.46%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.11'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.76%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.97'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.09%  '

# Row 45
$ws.Range("E45").Value = '  +2.86%  '

# Row 46
$ws.Range("E46").Value = '  +7.07%  '

# Row 47
$ws.Range("E47").Value = '  +0.02%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.40'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.02%  '

# Row 49
$ws.Range("E49").Value = '  -0.66%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '144.59'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.09%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.82'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.61%  '
